$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates (Nile Air NP-119 threat)
$ws.Range("A2").Value = "'29-JAN-26"
$ws.Range("D2").Value = 6022
$ws.Range("E2").Value = 6237
$ws.Range("F2").Value = -215

# Row 3 updates (Nile Air NP-109 threat)
$ws.Range("D3").Value = 8173
$ws.Range("E3").Value = 12107
$ws.Range("F3").Value = -3934
